# Weekly update: insert a new "Ajo" (garlic) price record for
# "Terminal La Palmera de La Serena" at row 114, pushing the existing
# records (old rows 114-145) down by one row (new rows 115-146).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 114 - everything that was
# at row 114 onwards shifts down to make room (dimension grows from
# A1:R145 to A1:R146).
$ws.Rows.Item(114).Insert()

# Fill the newly inserted row 114 with the new weekly data point.
$ws.Cells.Item(114, 1).Value = 8
$ws.Cells.Item(114, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(114, 3).Value = "Coquimbo"
$ws.Cells.Item(114, 4).Value = 44463
$ws.Cells.Item(114, 5).Value = 4
$ws.Cells.Item(114, 6).Value = 100112003
$ws.Cells.Item(114, 7).Value = "Ajo"
$ws.Cells.Item(114, 8).Value = "Chino"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 700
$ws.Cells.Item(114, 11).Value = 15000
$ws.Cells.Item(114, 12).Value = 16000
$ws.Cells.Item(114, 13).Value = 15500
$ws.Cells.Item(114, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(114, 15).Value = "China"
$ws.Cells.Item(114, 16).Value = 1550
$ws.Cells.Item(114, 17).Value = 10
$ws.Cells.Item(114, 18).Value = "Hortaliza"
